# More of the application skeleton added, including CLI and packaging.
# Currently works for tiling the basic "Standalone Template" project using
# wall-to-wall config file.
#
# Functional changes applied to docs/config_parsing_and_defaults.xlsx:
#   - Row 33 (bounding_box section header row) now starts a new section:
#       A33: "<root>"      -> "bounding_box"
#       B33: "bounding_box" -> "(layer def keys)"
#   - Column F ("Behaviour") narrowed from 67 chars to ~50.86 chars wide.
#   - Selection/scroll moved from the bottom of the sheet (A60) up near the
#     new bounding_box section (topLeftCell A10, active cell B33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (row 33) ---------------------------------------
$ws.Range("A33").Value = "bounding_box"
$ws.Range("B33").Value = "(layer def keys)"

# --- Column F width --------------------------------------------------------
# Target stored width is 50.85546875 characters. The Excel COM ColumnWidth
# setter quantizes to on-screen pixel boundaries, so we pick an input value
# that lands in the pixel bucket closest to the target width.
$ws.Columns.Item(6).ColumnWidth = 50.02

# --- Selection / scroll position -------------------------------------------
[void]$ws.Range("B33").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
